# Weekly data refresh: a new "Mora" price observation (week of 2022-01-21)
# is inserted at row 53 of the sheet, pushing the existing rows 53:73 down
# to 54:74 (dimension grows from A1:T73 to A1:T74).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 53, shifting rows 53-73 -> 54-74.
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new observation.
$ws.Range("A53").Value = 6
$ws.Range("B53").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C53").Value = "Metropolitana"
$ws.Range("D53").Value = 44582
$ws.Range("E53").Value = 13
$ws.Range("F53").Value = "Fruta"
$ws.Range("G53").Value = 100101
$ws.Range("H53").Value = "Berries"
$ws.Range("I53").Value = 100101008
$ws.Range("J53").Value = "Mora"
$ws.Range("K53").Value = "Sin especificar"
$ws.Range("L53").Value = "Primera"
$ws.Range("M53").Value = 200
$ws.Range("N53").Value = 6000
$ws.Range("O53").Value = 6000
$ws.Range("P53").Value = 6000
$ws.Range("Q53").Value = "$/bandeja 2 kilos"
$ws.Range("R53").Value = "Provincia de Curicó"
$ws.Range("S53").Value = 3000
$ws.Range("T53").Value = 2
